$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new blank worksheet column at H ("נהג") — this shifts the
# existing "נהג במקור"/"ממליץ"/"הערות" columns one position to the right
# (H->I, I->J, J->K), carrying their values/styles/column-widths along.
$ws.Columns("H").Insert()

# Grow the table to cover the new column (A1:J2 -> A1:K2).
$lo.Resize($ws.Range("A1:K2"))

# Re-stamp every header cell from H to K so the table's column metadata
# (names) picks up the right text for each (now shifted) position.
$ws.Range("H1").Value = "נהג"
$ws.Range("I1").Value = "נהג במקור"
$ws.Range("J1").Value = "ממליץ"
$ws.Range("K1").Value = "הערות"

# New driver column's sample data value for row 2.
$ws.Range("H2").Value = "פלוני"
# Re-stamp the shifted data cell too so the table metadata lines up.
$ws.Range("I2").Value = "נחום"

# Match the author's final selection.
[void]$ws.Range("H3").Select()
